$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty row 4 ("Survey 3") with its results
$ws.Range("B4").Value = 37
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = 19
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Move the active selection to F14 (matches the saved cursor position)
$ws.Range("F14").Select()
